# Updates cryptos list values (Price and Volume(1h) columns) to match
# the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): some values look numeric (e.g. "1.00", "342.00") and
# would otherwise be auto-converted/reformatted by Excel's type inference when
# assigned through .Value. Force the cell to text first, assign the exact
# string, then restore the default "Normal" style so no stray formatting is
# left behind (matches the original unstyled cells).
$ws.Range("D2").NumberFormat = "@"; $ws.Range("D2").Value = "63.390.00"; $ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"; $ws.Range("D3").Value = "2.659.12"; $ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"; $ws.Range("D5").Value = "610.44"; $ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"; $ws.Range("D6").Value = "143.87"; $ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"; $ws.Range("D9").Value = "2.657.74"; $ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"; $ws.Range("D11").Value = "5.63"; $ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"; $ws.Range("D14").Value = "27.37"; $ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"; $ws.Range("D15").Value = "3.133.56"; $ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"; $ws.Range("D16").Value = "63.228.75"; $ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"; $ws.Range("D17").Value = "0.0000145"; $ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"; $ws.Range("D18").Value = "2.661.22"; $ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"; $ws.Range("D20").Value = "342.24"; $ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"; $ws.Range("D23").Value = "0.999"; $ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"; $ws.Range("D24").Value = "66.95"; $ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"; $ws.Range("D25").Value = "1.65"; $ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"; $ws.Range("D27").Value = "8.66"; $ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"; $ws.Range("D29").Value = "548.06"; $ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"; $ws.Range("D30").Value = "1.00"; $ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"; $ws.Range("D31").Value = "7.83"; $ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"; $ws.Range("D32").Value = "2.05"; $ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"; $ws.Range("D34").Value = "0.0₃0808"; $ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"; $ws.Range("D35").Value = "173.00"; $ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"; $ws.Range("D36").Value = "5.14"; $ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"; $ws.Range("D37").Value = "0.406"; $ws.Range("D37").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"; $ws.Range("D41").Value = "174.75"; $ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"; $ws.Range("D42").Value = "0.999"; $ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"; $ws.Range("D43").Value = "3.75"; $ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"; $ws.Range("D44").Value = "22.22"; $ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"; $ws.Range("D45").Value = "0.0572"; $ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"; $ws.Range("D46").Value = "0.633"; $ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"; $ws.Range("D47").Value = "0.0962"; $ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"; $ws.Range("D49").Value = "18.75"; $ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"; $ws.Range("D50").Value = "1.76"; $ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"; $ws.Range("D51").Value = "11.27"; $ws.Range("D51").Style = "Normal"

# --- Volume(1h) column (E): values are plain text (contain '%' and padding
# spaces) so a direct assignment keeps them as text without extra handling.
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +4.30%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("E15").Value = "  +3.23%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("E19").Value = "  +3.77%  "
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +5.78%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  +16.76%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  +5.80%  "
$ws.Range("E33").Value = "  +7.26%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  +13.38%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  +9.33%  "
$ws.Range("E41").Value = "  +10.99%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("E44").Value = "  +5.00%  "
$ws.Range("E45").Value = "  +6.40%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("E51").Value = "  -1.01%  "
